$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.041.72'
$ws.Range('E2').Value = '  +0.18%  '
$ws.Range('D3').Value = '2.957.97'
$ws.Range('E3').Value = '  +1.00%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '380.60'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.03%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '102.55'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.59%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.543'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.80%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.586'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.67%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.45'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.17%  '
$ws.Range('E11').Value = '  -0.38%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0850'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.08%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '18.43'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.07%  '
$ws.Range('D14').Value = '3.419.75'
$ws.Range('E14').Value = '  +0.99%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '12.43'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +75.13%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '7.75'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +5.94%  '
$ws.Range('D17').Value = '2.957.86'
$ws.Range('E17').Value = '  -2.70%  '
$ws.Range('E18').Value = '  +4.41%  '
$ws.Range('D19').Value = '51.071.50'
$ws.Range('E19').Value = '  +0.37%  '
$ws.Range('E20').Value = '  -2.06%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.38'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.63%  '
$ws.Range('E22').Value = '  +0.79%  '
$ws.Range('E23').Value = '  +18.69%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '69.67'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.38%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '267.00'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.20%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.99'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.11%  '
$ws.Range('E27').Value = '  -0.02%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '25.84'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.38%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.165'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.86%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.95'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -6.98%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.108'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -4.35%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '10.50'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +7.54%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '50.73'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.17%  '
$ws.Range('B34').Value = 'InjectiveProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '33.97'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.68%  '
$ws.Range('B35').Value = 'Toncoin'
$ws.Range('C35').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.06'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.07%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0435'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.43%  '
$ws.Range('E37').Value = '  -0.02%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.20'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +8.22%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '16.72'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +3.11%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.117'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.04%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.83'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.03%  '
$ws.Range('E42').Value = '  -2.96%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '119.94'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.81%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.56'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +11.74%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '21.52'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.53%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.03'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.94%  '
$ws.Range('D47').Value = '2.028.26'
$ws.Range('E47').Value = '  +1.43%  '
$ws.Range('E48').Value = '  -1.62%  '
$ws.Range('E49').Value = '  -4.56%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0321'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -6.90%  '
$ws.Range('E51').Value = '  +7.40%  '
